$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear existing data rows (keep header row 1 intact) ---
$ws.Range("A2:T5").ClearContents()

# --- Re-introduce the cluster-name strings in the desired shared-string
#     order (ECs, FAPs, sCs, Spp1, Itga4) via a scratch area, then clear it. ---
$ws.Cells.Item(50,1).Value = "ECs"
$ws.Cells.Item(50,2).Value = "FAPs"
$ws.Cells.Item(50,3).Value = "sCs"
$ws.Cells.Item(50,4).Value = "Spp1"
$ws.Cells.Item(50,5).Value = "Itga4"
$ws.Range("A50:E50").ClearContents()

# --- Write the new data table (rows 2-10, columns A-T) ---
# Row 2
$ws.Cells.Item(2,"A").Value = "ECs"
$ws.Cells.Item(2,"B").Value = "Spp1"
$ws.Cells.Item(2,"C").Value = "Itga4"
$ws.Cells.Item(2,"D").Value = "ECs"
$ws.Cells.Item(2,"E").Value = 2
$ws.Cells.Item(2,"F").Value = 0.6666666666666666
$ws.Cells.Item(2,"G").Value = 209.0063303333334
$ws.Cells.Item(2,"H").Value = 627.018991
$ws.Cells.Item(2,"I").Value = 0.6751081226665357
$ws.Cells.Item(2,"J").Value = 0.6751081226665357
$ws.Cells.Item(2,"K").Value = 3
$ws.Cells.Item(2,"L").Value = 1
$ws.Cells.Item(2,"M").Value = 22.906497
$ws.Cells.Item(2,"N").Value = 68.719491
$ws.Cells.Item(2,"O").Value = 0.9446038650914245
$ws.Cells.Item(2,"P").Value = 0.9446038650914245
$ws.Cells.Item(2,"Q").Value = 4787.60287876151
$ws.Cells.Item(2,"R").Value = 43088.42590885358
$ws.Cells.Item(2,"S").Value = 0.6377097420254252
$ws.Cells.Item(2,"T").Value = 0.6377097420254252

# Row 3
$ws.Cells.Item(3,"A").Value = "ECs"
$ws.Cells.Item(3,"B").Value = "Spp1"
$ws.Cells.Item(3,"C").Value = "Itga4"
$ws.Cells.Item(3,"D").Value = "FAPs"
$ws.Cells.Item(3,"E").Value = 2
$ws.Cells.Item(3,"F").Value = 0.6666666666666666
$ws.Cells.Item(3,"G").Value = 209.0063303333334
$ws.Cells.Item(3,"H").Value = 627.018991
$ws.Cells.Item(3,"I").Value = 0.6751081226665357
$ws.Cells.Item(3,"J").Value = 0.6751081226665357
$ws.Cells.Item(3,"K").Value = 2
$ws.Cells.Item(3,"L").Value = 0.6666666666666666
$ws.Cells.Item(3,"M").Value = 0.1329193333333333
$ws.Cells.Item(3,"N").Value = 0.3987579999999999
$ws.Cells.Item(3,"O").Value = 0.005481244732096839
$ws.Cells.Item(3,"P").Value = 0.005481244732096839
$ws.Cells.Item(3,"Q").Value = 27.78098209035311
$ws.Cells.Item(3,"R").Value = 250.028838813178
$ws.Cells.Item(3,"S").Value = 0.003700432840961736
$ws.Cells.Item(3,"T").Value = 0.003700432840961736

# Row 4
$ws.Cells.Item(4,"A").Value = "ECs"
$ws.Cells.Item(4,"B").Value = "Spp1"
$ws.Cells.Item(4,"C").Value = "Itga4"
$ws.Cells.Item(4,"D").Value = "sCs"
$ws.Cells.Item(4,"E").Value = 2
$ws.Cells.Item(4,"F").Value = 0.6666666666666666
$ws.Cells.Item(4,"G").Value = 209.0063303333334
$ws.Cells.Item(4,"H").Value = 627.018991
$ws.Cells.Item(4,"I").Value = 0.6751081226665357
$ws.Cells.Item(4,"J").Value = 0.6751081226665357
$ws.Cells.Item(4,"K").Value = 3
$ws.Cells.Item(4,"L").Value = 1
$ws.Cells.Item(4,"M").Value = 1.210428333333333
$ws.Cells.Item(4,"N").Value = 3.631285
$ws.Cells.Item(4,"O").Value = 0.04991489017647865
$ws.Cells.Item(4,"P").Value = 0.04991489017647865
$ws.Cells.Item(4,"Q").Value = 252.9871840814928
$ws.Cells.Item(4,"R").Value = 2276.884656733435
$ws.Cells.Item(4,"S").Value = 0.03369794780014881
$ws.Cells.Item(4,"T").Value = 0.03369794780014881

# Row 5
$ws.Cells.Item(5,"A").Value = "FAPs"
$ws.Cells.Item(5,"B").Value = "Spp1"
$ws.Cells.Item(5,"C").Value = "Itga4"
$ws.Cells.Item(5,"D").Value = "ECs"
$ws.Cells.Item(5,"E").Value = 3
$ws.Cells.Item(5,"F").Value = 1
$ws.Cells.Item(5,"G").Value = 1.401741666666666
$ws.Cells.Item(5,"H").Value = 4.205225
$ws.Cells.Item(5,"I").Value = 0.004527744128790482
$ws.Cells.Item(5,"J").Value = 0.004527744128790482
$ws.Cells.Item(5,"K").Value = 3
$ws.Cells.Item(5,"L").Value = 1
$ws.Cells.Item(5,"M").Value = 22.906497
$ws.Cells.Item(5,"N").Value = 68.719491
$ws.Cells.Item(5,"O").Value = 0.9446038650914245
$ws.Cells.Item(5,"P").Value = 0.9446038650914245
$ws.Cells.Item(5,"Q").Value = 32.108991282275
$ws.Cells.Item(5,"R").Value = 288.980921540475
$ws.Cells.Item(5,"S").Value = 0.004276924604200494
$ws.Cells.Item(5,"T").Value = 0.004276924604200494

# Row 6
$ws.Cells.Item(6,"A").Value = "FAPs"
$ws.Cells.Item(6,"B").Value = "Spp1"
$ws.Cells.Item(6,"C").Value = "Itga4"
$ws.Cells.Item(6,"D").Value = "FAPs"
$ws.Cells.Item(6,"E").Value = 3
$ws.Cells.Item(6,"F").Value = 1
$ws.Cells.Item(6,"G").Value = 1.401741666666666
$ws.Cells.Item(6,"H").Value = 4.205225
$ws.Cells.Item(6,"I").Value = 0.004527744128790482
$ws.Cells.Item(6,"J").Value = 0.004527744128790482
$ws.Cells.Item(6,"K").Value = 2
$ws.Cells.Item(6,"L").Value = 0.6666666666666666
$ws.Cells.Item(6,"M").Value = 0.1329193333333333
$ws.Cells.Item(6,"N").Value = 0.3987579999999999
$ws.Cells.Item(6,"O").Value = 0.005481244732096839
$ws.Cells.Item(6,"P").Value = 0.005481244732096839
$ws.Cells.Item(6,"Q").Value = 0.1863185678388888
$ws.Cells.Item(6,"R").Value = 1.676867110549999
$ws.Cells.Item(6,"S").Value = 0.00002481767365421522
$ws.Cells.Item(6,"T").Value = 0.00002481767365421522

# Row 7
$ws.Cells.Item(7,"A").Value = "FAPs"
$ws.Cells.Item(7,"B").Value = "Spp1"
$ws.Cells.Item(7,"C").Value = "Itga4"
$ws.Cells.Item(7,"D").Value = "sCs"
$ws.Cells.Item(7,"E").Value = 3
$ws.Cells.Item(7,"F").Value = 1
$ws.Cells.Item(7,"G").Value = 1.401741666666666
$ws.Cells.Item(7,"H").Value = 4.205225
$ws.Cells.Item(7,"I").Value = 0.004527744128790482
$ws.Cells.Item(7,"J").Value = 0.004527744128790482
$ws.Cells.Item(7,"K").Value = 3
$ws.Cells.Item(7,"L").Value = 1
$ws.Cells.Item(7,"M").Value = 1.210428333333333
$ws.Cells.Item(7,"N").Value = 3.631285
$ws.Cells.Item(7,"O").Value = 0.04991489017647865
$ws.Cells.Item(7,"P").Value = 0.04991489017647865
$ws.Cells.Item(7,"Q").Value = 1.696707829347222
$ws.Cells.Item(7,"R").Value = 15.270370464125
$ws.Cells.Item(7,"S").Value = 0.0002260018509357729
$ws.Cells.Item(7,"T").Value = 0.0002260018509357729

# Row 8
$ws.Cells.Item(8,"A").Value = "sCs"
$ws.Cells.Item(8,"B").Value = "Spp1"
$ws.Cells.Item(8,"C").Value = "Itga4"
$ws.Cells.Item(8,"D").Value = "ECs"
$ws.Cells.Item(8,"E").Value = 3
$ws.Cells.Item(8,"F").Value = 1
$ws.Cells.Item(8,"G").Value = 99.18134533333334
$ws.Cells.Item(8,"H").Value = 297.544036
$ws.Cells.Item(8,"I").Value = 0.3203641332046738
$ws.Cells.Item(8,"J").Value = 0.3203641332046737
$ws.Cells.Item(8,"K").Value = 3
$ws.Cells.Item(8,"L").Value = 1
$ws.Cells.Item(8,"M").Value = 22.906497
$ws.Cells.Item(8,"N").Value = 68.719491
$ws.Cells.Item(8,"O").Value = 0.9446038650914245
$ws.Cells.Item(8,"P").Value = 0.9446038650914245
$ws.Cells.Item(8,"Q").Value = 2271.897189333964
$ws.Cells.Item(8,"R").Value = 20447.07470400568
$ws.Cells.Item(8,"S").Value = 0.3026171984617988
$ws.Cells.Item(8,"T").Value = 0.3026171984617987

# Row 9
$ws.Cells.Item(9,"A").Value = "sCs"
$ws.Cells.Item(9,"B").Value = "Spp1"
$ws.Cells.Item(9,"C").Value = "Itga4"
$ws.Cells.Item(9,"D").Value = "FAPs"
$ws.Cells.Item(9,"E").Value = 3
$ws.Cells.Item(9,"F").Value = 1
$ws.Cells.Item(9,"G").Value = 99.18134533333334
$ws.Cells.Item(9,"H").Value = 297.544036
$ws.Cells.Item(9,"I").Value = 0.3203641332046738
$ws.Cells.Item(9,"J").Value = 0.3203641332046737
$ws.Cells.Item(9,"K").Value = 2
$ws.Cells.Item(9,"L").Value = 0.6666666666666666
$ws.Cells.Item(9,"M").Value = 0.1329193333333333
$ws.Cells.Item(9,"N").Value = 0.3987579999999999
$ws.Cells.Item(9,"O").Value = 0.005481244732096839
$ws.Cells.Item(9,"P").Value = 0.005481244732096839
$ws.Cells.Item(9,"Q").Value = 13.18311830080978
$ws.Cells.Item(9,"R").Value = 118.648064707288
$ws.Cells.Item(9,"S").Value = 0.001755994217480888
$ws.Cells.Item(9,"T").Value = 0.001755994217480888

# Row 10
$ws.Cells.Item(10,"A").Value = "sCs"
$ws.Cells.Item(10,"B").Value = "Spp1"
$ws.Cells.Item(10,"C").Value = "Itga4"
$ws.Cells.Item(10,"D").Value = "sCs"
$ws.Cells.Item(10,"E").Value = 3
$ws.Cells.Item(10,"F").Value = 1
$ws.Cells.Item(10,"G").Value = 99.18134533333334
$ws.Cells.Item(10,"H").Value = 297.544036
$ws.Cells.Item(10,"I").Value = 0.3203641332046738
$ws.Cells.Item(10,"J").Value = 0.3203641332046737
$ws.Cells.Item(10,"K").Value = 3
$ws.Cells.Item(10,"L").Value = 1
$ws.Cells.Item(10,"M").Value = 1.210428333333333
$ws.Cells.Item(10,"N").Value = 3.631285
$ws.Cells.Item(10,"O").Value = 0.04991489017647865
$ws.Cells.Item(10,"P").Value = 0.04991489017647865
$ws.Cells.Item(10,"Q").Value = 120.0519105295845
$ws.Cells.Item(10,"R").Value = 1080.46719476626
$ws.Cells.Item(10,"S").Value = 0.01599094052539407
$ws.Cells.Item(10,"T").Value = 0.01599094052539406
